{"js": "// Replace each two-digit-division problem's text in the document with its\n// new value. Every source string is unique in the document, so a simple\n// body.search()+insertText(\"Replace\") pass for each pair is safe.\nconst pairs = [\n  [\"25\u00f78=\", \"98\u00f75=\"],\n  [\"67\u00f77=\", \"37\u00f75=\"],\n  [\"44\u00f75=\", \"62\u00f72=\"],\n  [\"11\u00f72=\", \"30\u00f74=\"],\n  [\"41\u00f77=\", \"88\u00f74=\"],\n  [\"78\u00f73=\", \"67\u00f75=\"],\n  [\"69\u00f77=\", \"98\u00f77=\"],\n  [\"28\u00f79=\", \"86\u00f75=\"],\n  [\"60\u00f76=\", \"69\u00f75=\"],\n  [\"27\u00f75=\", \"13\u00f75=\"],\n  [\"97\u00f73=\", \"86\u00f72=\"],\n  [\"39\u00f79=\", \"80\u00f74=\"],\n  [\"86\u00f74=\", \"40\u00f73=\"],\n  [\"34\u00f72=\", \"26\u00f72=\"],\n  [\"36\u00f79=\", \"93\u00f78=\"],\n  [\"58\u00f76=\", \"25\u00f76=\"],\n  [\"55\u00f78=\", \"60\u00f74=\"],\n  [\"40\u00f75=\", \"34\u00f77=\"],\n  [\"27\u00f77=\", \"99\u00f79=\"],\n  [\"21\u00f79=\", \"65\u00f78=\"],\n  [\"77\u00f74=\", \"64\u00f76=\"],\n  [\"81\u00f74=\", \"54\u00f74=\"],\n  [\"85\u00f73=\", \"24\u00f75=\"],\n  [\"91\u00f77=\", \"69\u00f75=\"],\n  [\"44\u00f77=\", \"79\u00f78=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text \"${oldText}\" to replace.`);\n  }\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-division problem's text in the document with its\n# new value. Every source string is unique in the document, so a simple\n# Find/Replace pass per pair is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"25\u00f78=\", \"98\u00f75=\"),\n    @(\"67\u00f77=\", \"37\u00f75=\"),\n    @(\"44\u00f75=\", \"62\u00f72=\"),\n    @(\"11\u00f72=\", \"30\u00f74=\"),\n    @(\"41\u00f77=\", \"88\u00f74=\"),\n    @(\"78\u00f73=\", \"67\u00f75=\"),\n    @(\"69\u00f77=\", \"98\u00f77=\"),\n    @(\"28\u00f79=\", \"86\u00f75=\"),\n    @(\"60\u00f76=\", \"69\u00f75=\"),\n    @(\"27\u00f75=\", \"13\u00f75=\"),\n    @(\"97\u00f73=\", \"86\u00f72=\"),\n    @(\"39\u00f79=\", \"80\u00f74=\"),\n    @(\"86\u00f74=\", \"40\u00f73=\"),\n    @(\"34\u00f72=\", \"26\u00f72=\"),\n    @(\"36\u00f79=\", \"93\u00f78=\"),\n    @(\"58\u00f76=\", \"25\u00f76=\"),\n    @(\"55\u00f78=\", \"60\u00f74=\"),\n    @(\"40\u00f75=\", \"34\u00f77=\"),\n    @(\"27\u00f77=\", \"99\u00f79=\"),\n    @(\"21\u00f79=\", \"65\u00f78=\"),\n    @(\"77\u00f74=\", \"64\u00f76=\"),\n    @(\"81\u00f74=\", \"54\u00f74=\"),\n    @(\"85\u00f73=\", \"24\u00f75=\"),\n    @(\"91\u00f77=\", \"69\u00f75=\"),\n    @(\"44\u00f77=\", \"79\u00f78=\")\n)\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1            # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute([ref]$findText, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$replaceText, [ref]2) | Out-Null\n}\n"}
